$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-9 down to 6-10.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 (same Mercado/Region/Producto metadata as the other rows).
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 45002
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107011
$ws.Range("J5").Value = "Tuna"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 12500
$ws.Range("Q5").Value = "`$/caja 18 kilos"
$ws.Range("R5").Value = "Provincia de Melipilla"
$ws.Range("S5").Value = 694
$ws.Range("T5").Value = 18
